$d = $word.ActiveDocument

# --- Simple text replacements (visible content only changes) ---

# "Modulo para atender llamadas" -> "Módulo para atender llamadas"
$d.Content.Find.Execute("Modulo para atender llamadas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Módulo para atender llamadas", 2)

# "Necesitábamos un modulo que" -> "Necesitábamos un módulo que"
$d.Content.Find.Execute("Necesitábamos un modulo que", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Necesitábamos un módulo que", 2)

# Status: "Pendiente" -> "Aceptada"
$d.Content.Find.Execute("Pendiente", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Aceptada", 2)

# --- Fill previously-empty cells in table 1 ---
$t = $d.Tables.Item(1)

# Row 9: Alternative decisions (options) -> "-"
$t.Cell(9, 2).Range.Text = "-"

# Row 10: Decision outcome (options selected) -> "ADD-041"
$t.Cell(10, 2).Range.Text = "ADD-041"

# Row 11: Pros opciones -> "Versátil, Interfaz amistosa"
$t.Cell(11, 2).Range.Text = "Versátil, Interfaz amistosa"

# Row 12: Cons opciones -> "Solo dispositivos Android"
$t.Cell(12, 2).Range.Text = "Solo dispositivos Android"
